# Updated cryptos list with refreshed prices and 1h volume percentages.
# D-column values are prefixed with a leading apostrophe so Excel stores
# them as text (matching the original text-formatted "Price" column)
# instead of auto-converting them to numbers and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.038.26"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'2.361.36"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.680"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'239.58"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.595"
$ws.Range("E9").Value = "  +10.19%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").Value = "'57.21"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'32.23"
$ws.Range("E12").Value = "  +13.13%  "
$ws.Range("D13").Value = "'7.30"
$ws.Range("E13").Value = "  +9.91%  "
$ws.Range("D15").Value = "'2.710.35"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "'16.67"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'2.366.39"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'43.890.91"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +5.33%  "
$ws.Range("D22").Value = "'77.04"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "'256.77"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  +24.80%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'10.73"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "'175.06"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").Value = "'0.0761"
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "'5.43"
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("E37").Value = "  -7.64%  "
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'0.0278"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").Value = "'0.113"
$ws.Range("E41").Value = "  +15.34%  "
$ws.Range("D42").Value = "'0.206"
$ws.Range("E42").Value = "  +13.98%  "
$ws.Range("D43").Value = "'9.17"
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("D44").Value = "'19.19"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "'4.74"
$ws.Range("E46").Value = "  +6.50%  "
$ws.Range("D47").Value = "'58.72"
$ws.Range("E47").Value = "  +11.75%  "
$ws.Range("E48").Value = "  +7.80%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").Value = "'100.43"
$ws.Range("E51").Value = "  +2.64%  "
